# Update "想去人数" (F column) figures across the 展览, 演出 and 全部类型 sheets
# to reflect newly generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsPerformance = $wb.Worksheets.Item("演出")
$wsAll = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibition) sheet updates: row -> new F value
$exhibitionUpdates = @{
    2  = 1122
    4  = 256
    5  = 1794
    6  = 673
    7  = 526
    8  = 4829
    9  = 58
    12 = 990
    13 = 329
    14 = 1291
    17 = 2996
    18 = 1830
    22 = 60
    24 = 944
    27 = 2979
    28 = 1030
    29 = 2505
    31 = 1344
    32 = 3640
    33 = 94
    34 = 896
    36 = 1145
    38 = 1197
    39 = 25
    40 = 887
    41 = 567
    42 = 278
    43 = 371
    44 = 290
    45 = 3499
}

foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# 演出 (Performance) sheet updates
$performanceUpdates = @{
    3 = 16
}

foreach ($row in $performanceUpdates.Keys) {
    $wsPerformance.Range("F$row").Value = $performanceUpdates[$row]
}

# 全部类型 (All types) sheet updates
$allTypesUpdates = @{
    2  = 1122
    4  = 256
    6  = 1794
    7  = 673
    8  = 526
    9  = 4829
    10 = 58
    11 = 16
    13 = 329
    14 = 1291
    15 = 2996
    17 = 1830
    25 = 60
    26 = 944
    28 = 2979
    30 = 1030
    31 = 2505
    32 = 1344
    33 = 3640
    35 = 94
    36 = 896
    37 = 1145
    40 = 1197
    41 = 887
    42 = 567
    43 = 371
    47 = 290
    48 = 3499
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allTypesUpdates[$row]
}
